$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date-formatted style (s="4") from the row above into the new rows' H:I columns
$ws.Range("H230:I230").Copy()
$ws.Range("H231:I235").PasteSpecial(-4122)

# Row 231: 370. Range Addition
$ws.Cells.Item(231, 1).Value = 370
$ws.Cells.Item(231, 2).Value = "Range Addition"
$ws.Cells.Item(231, 3).Value = "#array #prefix-sum "
$ws.Cells.Item(231, 4).Value = "medium"
$ws.Cells.Item(231, 5).Value = 0
$ws.Cells.Item(231, 6).Value = 1
$ws.Cells.Item(231, 7).Value = 30
$ws.Cells.Item(231, 8).Value = 45990
$ws.Cells.Item(231, 9).Value = 45990

# Row 232: 1590. Make Sum Divisible by P
$ws.Cells.Item(232, 1).Value = 1590
$ws.Cells.Item(232, 2).Value = "Make Sum Divisible by P"
$ws.Cells.Item(232, 3).Value = "#array #prefix-sum #divide "
$ws.Cells.Item(232, 4).Value = "medium"
$ws.Cells.Item(232, 5).Value = 0
$ws.Cells.Item(232, 6).Value = 1
$ws.Cells.Item(232, 7).Value = 25
$ws.Cells.Item(232, 8).Value = 45991
$ws.Cells.Item(232, 9).Value = 45991

# Row 233: 2141. Maximum Running Time of N Computers
$ws.Cells.Item(233, 1).Value = 2141
$ws.Cells.Item(233, 2).Value = "Maximum Running Time of N Computers"
$ws.Cells.Item(233, 3).Value = "#array #binary-search #greedy "
$ws.Cells.Item(233, 4).Value = "hard"
$ws.Cells.Item(233, 5).Value = 0
$ws.Cells.Item(233, 6).Value = 1
$ws.Cells.Item(233, 7).Value = "???"
$ws.Cells.Item(233, 8).Value = 45992
$ws.Cells.Item(233, 9).Value = 45992

# Row 234: 3623. Count Number of Trapezoids I
$ws.Cells.Item(234, 1).Value = 3623
$ws.Cells.Item(234, 2).Value = "Count Number of Trapezoids I"
$ws.Cells.Item(234, 4).Value = "medium"
$ws.Cells.Item(234, 5).Value = 0
$ws.Cells.Item(234, 6).Value = 1
$ws.Cells.Item(234, 7).Value = 45
$ws.Cells.Item(234, 8).Value = 45993
$ws.Cells.Item(234, 9).Value = 45993

# Row 235: 1214. Two Sum BSTs
$ws.Cells.Item(235, 1).Value = 1214
$ws.Cells.Item(235, 2).Value = "Two Sum BSTs"
$ws.Cells.Item(235, 3).Value = "#binary-tree #bst #bfs #dfs #morris"
$ws.Cells.Item(235, 4).Value = "medium"
$ws.Cells.Item(235, 5).Value = 1
$ws.Cells.Item(235, 6).Value = 0
$ws.Cells.Item(235, 7).Value = 16
$ws.Cells.Item(235, 8).Value = 45993
$ws.Cells.Item(235, 9).Value = 45993

$ws.Rows.Item(231).RowHeight = 17
$ws.Rows.Item(232).RowHeight = 34
$ws.Rows.Item(233).RowHeight = 34
$ws.Rows.Item(234).RowHeight = 34
$ws.Rows.Item(235).RowHeight = 34

$ws.Range("H235:I235").Select()
